$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8611539602279663
$ws.Range("B1").Value = 1.224928498268127
$ws.Range("C1").Value = 2.016730308532715
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 1.746039152145386
